# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 96
$ws1.Range("F4").Value = 1436
$ws1.Range("F5").Value = 178
$ws1.Range("F7").Value = 35
$ws1.Range("F8").Value = 9650
$ws1.Range("F9").Value = 162
$ws1.Range("F11").Value = 240
$ws1.Range("F14").Value = 6649
$ws1.Range("F15").Value = 1080
$ws1.Range("F16").Value = 121

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 96
$ws4.Range("F4").Value = 1436
$ws4.Range("F5").Value = 178
$ws4.Range("F7").Value = 35
$ws4.Range("F10").Value = 9650
$ws4.Range("F11").Value = 162
$ws4.Range("F13").Value = 240
$ws4.Range("F16").Value = 6649
$ws4.Range("F17").Value = 1080
$ws4.Range("F18").Value = 121
